$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2014652014652015
$ws.Range("C2").Value = 0.5824175824175825
$ws.Range("J2").Value = 0.003663003663003663
$ws.Range("P2").Value = 0.1282051282051282
$ws.Range("S2").Value = 0.08424908424908426
$ws.Range("B3").Value = 0.006060606060606061
$ws.Range("C3").Value = 0.02424242424242424
$ws.Range("J3").Value = 0.006060606060606061
$ws.Range("P3").Value = 0.8424242424242424
$ws.Range("S3").Value = 0.1212121212121212
$ws.Range("J4").Value = 0.04
$ws.Range("P4").Value = 0.72
$ws.Range("S4").Value = 0.24
$ws.Range("B6").Value = 0.05825242718446602
$ws.Range("D6").Value = 0.01456310679611651
$ws.Range("F6").Value = 0.06310679611650485
$ws.Range("J6").Value = 0.2572815533980582
$ws.Range("O6").Value = 0.02912621359223301
$ws.Range("Q6").Value = 0.1747572815533981
$ws.Range("R6").Value = 0.1116504854368932
$ws.Range("S6").Value = 0.2912621359223301
$ws.Range("B7").Value = 0.1037735849056604
$ws.Range("D7").Value = 0.02358490566037736
$ws.Range("E7").Value = 0.004716981132075472
$ws.Range("F7").Value = 0.05188679245283019
$ws.Range("J7").Value = 0.1556603773584906
$ws.Range("O7").Value = 0.009433962264150943
$ws.Range("Q7").Value = 0.1320754716981132
$ws.Range("R7").Value = 0.07547169811320754
$ws.Range("S7").Value = 0.4433962264150944
$ws.Range("B8").Value = 0.06403940886699508
$ws.Range("D8").Value = 0.01970443349753695
$ws.Range("E8").Value = 0.002463054187192118
$ws.Range("F8").Value = 0.05665024630541872
$ws.Range("J8").Value = 0.1477832512315271
$ws.Range("O8").Value = 0.01477832512315271
$ws.Range("Q8").Value = 0.1995073891625616
$ws.Range("R8").Value = 0.09359605911330049
$ws.Range("S8").Value = 0.4014778325123153
$ws.Range("B9").Value = 0.1764705882352941
$ws.Range("D9").Value = 0.0213903743315508
$ws.Range("F9").Value = 0.0374331550802139
$ws.Range("J9").Value = 0.1229946524064171
$ws.Range("O9").Value = 0.0106951871657754
$ws.Range("Q9").Value = 0.1122994652406417
$ws.Range("R9").Value = 0.1176470588235294
$ws.Range("S9").Value = 0.4010695187165775
$ws.Range("B10").Value = 0.09184423218221896
$ws.Range("D10").Value = 0.0227773695811903
$ws.Range("E10").Value = 0.001469507714915503
$ws.Range("F10").Value = 0.06171932402645114
$ws.Range("J10").Value = 0.1102130786186628
$ws.Range("O10").Value = 0.01322556943423953
$ws.Range("Q10").Value = 0.2167523879500367
$ws.Range("R10").Value = 0.106539309331374
$ws.Range("S10").Value = 0.3754592211609111
$ws.Range("G11").Value = 0.1470588235294118
$ws.Range("J11").Value = 0.07058823529411765
$ws.Range("K11").Value = 0.2147058823529412
$ws.Range("L11").Value = 0.55
$ws.Range("S11").Value = 0.01764705882352941
$ws.Range("G12").Value = 0.6974358974358974
$ws.Range("J12").Value = 0.2205128205128205
$ws.Range("K12").Value = 0.01538461538461539
$ws.Range("L12").Value = 0.03589743589743589
$ws.Range("S12").Value = 0.03076923076923077
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.2916666666666667
$ws.Range("S13").Value = 0.08333333333333333
$ws.Range("F15").Value = 0.02212389380530973
$ws.Range("H15").Value = 0.1858407079646018
$ws.Range("I15").Value = 0.04424778761061947
$ws.Range("J15").Value = 0.3584070796460177
$ws.Range("K15").Value = 0.08849557522123894
$ws.Range("M15").Value = 0.008849557522123894
$ws.Range("N15").Value = 0.004424778761061947
$ws.Range("O15").Value = 0.06637168141592921
$ws.Range("S15").Value = 0.2212389380530974
$ws.Range("H16").Value = 0.1568627450980392
$ws.Range("I16").Value = 0.09313725490196079
$ws.Range("J16").Value = 0.4215686274509804
$ws.Range("K16").Value = 0.1274509803921569
$ws.Range("M16").Value = 0.004901960784313725
$ws.Range("O16").Value = 0.03431372549019608
$ws.Range("S16").Value = 0.1617647058823529
$ws.Range("F17").Value = 0.0175054704595186
$ws.Range("H17").Value = 0.1597374179431072
$ws.Range("I17").Value = 0.0787746170678337
$ws.Range("J17").Value = 0.3982494529540481
$ws.Range("K17").Value = 0.1050328227571116
$ws.Range("M17").Value = 0.03719912472647702
$ws.Range("N17").Value = 0.002188183807439825
$ws.Range("O17").Value = 0.0700218818380744
$ws.Range("S17").Value = 0.1312910284463895
$ws.Range("F18").Value = 0.0163265306122449
$ws.Range("H18").Value = 0.1591836734693877
$ws.Range("I18").Value = 0.1142857142857143
$ws.Range("J18").Value = 0.4448979591836735
$ws.Range("K18").Value = 0.08979591836734693
$ws.Range("M18").Value = 0.004081632653061225
$ws.Range("O18").Value = 0.08163265306122448
$ws.Range("S18").Value = 0.08979591836734693
$ws.Range("F19").Value = 0.01576044129235618
$ws.Range("H19").Value = 0.1804570527974783
$ws.Range("I19").Value = 0.07249802994483845
$ws.Range("J19").Value = 0.4018912529550828
$ws.Range("K19").Value = 0.1150512214342002
$ws.Range("M19").Value = 0.02127659574468085
$ws.Range("N19").Value = 0.0007880220646178094
$ws.Range("O19").Value = 0.07407407407407407
$ws.Range("S19").Value = 0.1182033096926714
